# Zera os valores da planilha "Valores" (bug do excel apresentacao nao
# ser encontrado/preenchido ao gerar os relatorios por competencia/ano).
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Valores")

# Celulas que possuem valores diferentes de zero no relatorio original
# devem ser zeradas.
$cells = @("C1", "D1", "E1", "F1", "J1", "K1",
           "C4", "D4", "E4", "J4",
           "C5",
           "C6", "D6", "F6", "J6", "K6",
           "C16", "D16", "E16", "F16", "J16", "K16")

foreach ($cell in $cells) {
    $ws.Range($cell).Value = 0
}
